# Regenerate merged AHB files:
#  - rename the "_old" / "_new" column-header suffixes to the release
#    identifiers "_FV2210" / "_FV2304"
#  - freeze the header row
#  - turn the used range into a real Excel Table ("Table1") with an
#    autofilter on the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) rename the header row (A1:U1) -------------------------------------
$renameMap = @{
    "Segmentname_old"          = "Segmentname_FV2210"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2210"
    "Segment_old"              = "Segment_FV2210"
    "Datenelement_old"         = "Datenelement_FV2210"
    "Segment ID_old"           = "Segment ID_FV2210"
    "Code_old"                 = "Code_FV2210"
    "Qualifier_old"            = "Qualifier_FV2210"
    "Beschreibung_old"         = "Beschreibung_FV2210"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2210"
    "Bedingung_old"            = "Bedingung_FV2210"
    "Segmentname_new"          = "Segmentname_FV2304"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2304"
    "Segment_new"              = "Segment_FV2304"
    "Datenelement_new"         = "Datenelement_FV2304"
    "Segment ID_new"           = "Segment ID_FV2304"
    "Code_new"                 = "Code_FV2304"
    "Qualifier_new"            = "Qualifier_FV2304"
    "Beschreibung_new"         = "Beschreibung_FV2304"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2304"
    "Bedingung_new"            = "Bedingung_FV2304"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = $cell.Value2
    if ($renameMap.ContainsKey($cur)) {
        $cell.Value2 = $renameMap[$cur]
    }
}

# --- 2) freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) turn A1:U82 into a real table with an autofilter -------------------
$tableRange = $ws.Range("A1:U82")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
